$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert two new columns: new B ("Salt") and new D ("Pepper") ---
$ws.Columns("B").Insert()
$ws.Columns("D").Insert()

# Remove stray blank styled cells created by the column inserts
$ws.Range("B1:B7").Clear()
$ws.Range("D1:D7").Clear()

# --- GRADUAL HASHING block: fill in data rows 7-12 (columns C, E-J) ---
$ws.Range("C7").Value = "Testus"

$ws.Range("E7").Value = "b9ece18c950afbfa6b0fdbfa4ff731d3"
$ws.Range("E8").Value = "2408730ad248ad4e4aa36fb14f5e0631"
$ws.Range("E9").Value = "78ac12266411beaf4149f5cfb6e33afd"
$ws.Range("E10").Value = "0cbc6611f5540bd0809a388dc95a615b"
$ws.Range("E11").Value = "4e73c38e86323e3f00a309e6bcdcea46"
$ws.Range("E12").Value = "021034f85e1577648504c45f70486e2c"

$ws.Range("F7").Value = "c2c53d66948214258a26ca9ca845d7ac0c17f8e7"
$ws.Range("F8").Value = "1b1be06f5bfb451c3199aeb9e8e337905dd0cecb"
$ws.Range("F9").Value = "c82524e086d9a253cce031af2342715692b79aac"
$ws.Range("F10").Value = "640ab2bae07bedc4c163f679a746f7ab7fb5d1fa"
$ws.Range("F11").Value = "2d3d360a82c6464020c5cd02b2352b74c8b6df6f"
$ws.Range("F12").Value = "9d4b621d6602beb27e1fb82dbeb841dd7897661d"

$ws.Range("G7").Value = "e632b7095b0bf32c260fa4c539e9fd7b852d0de454e9be26f24d0d6f91d069d3"
$ws.Range("G8").Value = "7566ad8da86a586045385d1bdc1f1f26df9257280bf47a5c2975ee926dee3e6e"
$ws.Range("G9").Value = "47e8f094e1aca27195e6de64d47e8aa2b4f9c234d95da517084e0926fdaa8ba1"
$ws.Range("G10").Value = "532eaabd9574880dbf76b9b8cc00832c20a6ec113d682299550d7a6e0f345e25"
$ws.Range("G11").Value = "e835808c8a4b7b3e5c1fa276b121c382f74755a7c0dae908867209ef67fbce8b"
$ws.Range("G12").Value = "ca21ae9392cb39a8eddd8d76a33735e779a909aada3b66ccfa71fa8a22052f38"

$ws.Range("H7").Value = "b2396a002fe7aec008808687d7cbacb340b7f7a090008382f3c95870f6fb10415f61f5737c102d4bfec58fe525407ea2001e761dab1da8a501d9523921f0ec21"
$ws.Range("H8").Value = "cc12409c9f6cb05f7b8b959fcf5454ba457fccbf9fc83b278850a524b4962ee1ff81810f7ddb86fd16601c22548e3809dc440afd2701702268fcb91907fcff53"
$ws.Range("H9").Value = "71a0ec8b3318b49a9f5a12205a9cdedc9ac6b01a447c8b21c824aa33da564f1864a65476f03d1de71d12ea294555c810f9f9da8326c4ad13f151aa3ba7a33fd5"
$ws.Range("H10").Value = "c6ee9e33cf5c6715a1d148fd73f7318884b41adcb916021e2bc0e800a5c5dd97f5142178f6ae88c8fdd98e1afb0ce4c8d2c54b5f37b30b7da1997bb33b0b8a31"
$ws.Range("H11").Value = "726a250d106f6636303b67ee16fc2852bc28fb7c61cdc8b8c8da349a263fd2384cdb6e72e08cf62aa573b6a205483806180b57c3430b90bb8819738994e9306c"
$ws.Range("H12").Value = "6db16199111b65d307e5af1575a0e97dfecde1705a056dc92feb6e0b29444241996556d48359e7c59a369158462d4d0b151545667f701d365c45914cd8029f3c"

$ws.Range("I7").Value = "9de1976650ec90e603d2892e43015af720d91fe4"
$ws.Range("I8").Value = "77503acf53f822816d66efcf47be6a05027b495f"
$ws.Range("I9").Value = "dcbc4a87d4df24b0b640a91ffa15b83888ffe492"
$ws.Range("I10").Value = "76c82682cd7af7e812e513fa0e7914ab40b842e0"
$ws.Range("I11").Value = "a3b588877caf94021b02759a35b0684b30235221"
$ws.Range("I12").Value = "1460cbcc6d2447d17044b8d003d0b669fe999365"

$ws.Range("J7").Value = "be047a60"
$ws.Range("J8").Value = "a2d61f78"
$ws.Range("J9").Value = "4572e01a"
$ws.Range("J10").Value = "784dd132"
$ws.Range("J11").Value = "3ac2766f"
$ws.Range("J12").Value = "c63971b4"

# --- New header cells: Salt (B1) / Pepper (D1), bold font ---
$ws.Range("B1").Value = "Salt"
$ws.Range("B1").Font.Bold = $true
$ws.Range("B1").Font.Italic = $false

$ws.Range("D1").Value = "Pepper"
$ws.Range("D1").Font.Bold = $true
$ws.Range("D1").Font.Italic = $false

# --- HASHER SALT PEPPER block (rows 14-16) ---
$ws.Range("A14").Value = "HASHER SALT PEPPER"
$ws.Range("A14").Font.Bold = $true
$ws.Range("A14").Font.Italic = $true
$ws.Range("A14").Font.Size = 12

$ws.Range("C14").Value = "Data"
$ws.Range("B14").Value = "salt101"
$ws.Range("D15").Value = "M5"

$ws.Range("E14").Value = "b9f2d25febd39790a6a6b3163c672ec6"
$ws.Range("F14").Value = "a1ccf36317c6cf2d98bc03832ee35499e3dc62fe"
$ws.Range("G14").Value = "ad853ef97f34098a94f3f1ab59c6301c334ccc7b0c1ebfc2f6c5eb6bff60078d"
$ws.Range("H14").Value = "75a8914a2a4aeabea4cf8e1fef4d255f63708c96f118543cd0ef54db056d7b4f1ea9df4d0365119866c8ab8a85fcd45155b8108fb9c48e90074dbc8f5d19dd3a"
$ws.Range("J14").Value = "1a85bcb6"
$ws.Range("I14").Value = "7cef35065a9df17b8222781a616052cf6199c1b9"

$ws.Range("E15").Value = "d4a38bab94c4341032fd5109452635dc"
$ws.Range("F15").Value = "e9cdf0941adf6a4c11b29fec31b63210e2d79b29"
$ws.Range("G15").Value = "b9e79e9e860692da02b7a6e3cc0fd30f8d0bdbe185c0dbc5cdcd78e8514e176b"
$ws.Range("H15").Value = "3c3bf36c33157d048629df3483be44853d14514989db0c6d0b8eb38a73f927833b01a55d48a40806f2c7c628bf88347cc2c99b778154c29b318973793b4088a8"
$ws.Range("J15").Value = "3f75d493"
$ws.Range("I15").Value = "4e4aa21551256fddc151d31719cec0bade0f5b9c"

$ws.Range("E16").Value = "ffd60eb4265b259bd41cf1ecf9168d65"
$ws.Range("F16").Value = "610966b737670c4957cefa62884cbda5d06b6e94"
$ws.Range("G16").Value = "0967200e8d4e27894d2bf53d318a3efe2efada1de1c13dc910ef9558592e3858"
$ws.Range("H16").Value = "56af63202c76821d777a1aa088444098e8ff87f18cc36efdda5803cc7331ac0ed1c3d7daaf4c5465620f625d3f5a31c263c07b9c6bc15c9e48c8bea83558ac59"
$ws.Range("J16").Value = "22f50b88"
$ws.Range("I16").Value = "001bc637637344af6f63ca9ecf2cfe0fdbff8ba4"

# rows 15/16 also repeat "Data"/"salt101"/"M5" in columns C,D,B
$ws.Range("C15").Value = "Data"
$ws.Range("B16").Value = "salt101"
$ws.Range("C16").Value = "Data"
$ws.Range("D16").Value = "M5"

# --- page setup ---
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PaperSize = 9

# --- selection ---
$ws.Range("J16").Select() | Out-Null
